$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" status text to "In Translation" on all sheets
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        if ($cell.Value2 -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# The status column got narrower because the new text ("In Translation") is
# shorter than the old text ("Ready for handoff"), so Excel's column autofit
# shrank the columns that display it. Reproduce the resulting column widths.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E1").ColumnWidth = 12.5
$overview.Range("F1").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C1").ColumnWidth = 12.5
